$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 7951579.5
$ws.Range("J43").Value = 11111611
$ws.Range("L43").Value = 11111611
$ws.Range("N43").Value = -11111749
$ws.Range("H113").Value = 2688
$ws.Range("I113").Value = 2502.5
$ws.Range("J113").Value = 2762.2
$ws.Range("K113").Value = 2502.5
$ws.Range("L113").Value = 2762.2
$ws.Range("M113").Value = 751.5
$ws.Range("N113").Value = -9270.200000000001
$ws.Range("H129").Value = 894.6391599999999
$ws.Range("J129").Value = 945.43823
$ws.Range("L129").Value = 2836.31469
$ws.Range("N129").Value = -12836.31469
$ws.Range("H138").Value = 1410.02
$ws.Range("I138").Value = 602.9524
$ws.Range("J138").Value = 1994.4482
$ws.Range("K138").Value = 1808.8572
$ws.Range("L138").Value = 5983.3446
$ws.Range("M138").Value = 3331.1428
$ws.Range("N138").Value = -16263.3446
$ws.Range("H141").Value = 793.4286
$ws.Range("I141").Value = 793.4286
$ws.Range("K141").Value = 2380.2858
$ws.Range("M141").Value = 2799.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1144.6666
$ws.Range("I74").Value = 732.6
$ws.Range("J74").Value = 1968.8
$ws.Range("K74").Value = 732.6
$ws.Range("L74").Value = 1968.8
$ws.Range("M74").Value = 141.4
$ws.Range("N74").Value = -3716.8
$ws.Range("H77").Value = 1144.6666
$ws.Range("I77").Value = 732.6
$ws.Range("J77").Value = 1968.8
$ws.Range("K77").Value = 3663
$ws.Range("L77").Value = 9844
$ws.Range("M77").Value = 705
$ws.Range("N77").Value = -18580
$ws.Range("H135").Value = 20228.4
$ws.Range("J135").Value = 20228.4
$ws.Range("L135").Value = 20228.4
$ws.Range("N135").Value = -30368.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H20").Value = 1384.3125
$ws.Range("I20").Value = 1131.625
$ws.Range("K20").Value = 1131.625
$ws.Range("M20").Value = -884.625
$ws.Range("H86").Value = 3988.889
$ws.Range("I86").Value = 4030.4783
$ws.Range("K86").Value = 4030.4783
$ws.Range("M86").Value = -2907.4783
$ws.Range("H89").Value = 3988.889
$ws.Range("I89").Value = 4030.4783
$ws.Range("K89").Value = 20152.3915
$ws.Range("M89").Value = -14536.3915
$ws.Range("H107").Value = 1935.1538
$ws.Range("I107").Value = 1843
$ws.Range("J107").Value = 2042.6666
$ws.Range("K107").Value = 1843
$ws.Range("L107").Value = 2042.6666
$ws.Range("M107").Value = 77
$ws.Range("N107").Value = -5882.6666
$ws.Range("H134").Value = 2833.7542
$ws.Range("I134").Value = 818.7347
$ws.Range("J134").Value = 11061.75
$ws.Range("K134").Value = 2456.2041
$ws.Range("L134").Value = 33185.25
$ws.Range("M134").Value = 78.79590000000007
$ws.Range("N134").Value = -38255.25
$ws.Range("H138").Value = 64744.5
$ws.Range("J138").Value = 64744.5
$ws.Range("L138").Value = 64744.5
$ws.Range("N138").Value = -75024.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1359.2
$ws.Range("J58").Value = 1700
$ws.Range("L58").Value = 1700
$ws.Range("N58").Value = -2106
$ws.Range("H132").Value = 3122
$ws.Range("I132").Value = 1550
$ws.Range("J132").Value = 3571.1428
$ws.Range("K132").Value = 4650
$ws.Range("L132").Value = 10713.4284
$ws.Range("M132").Value = -2120
$ws.Range("N132").Value = -15773.4284
$ws.Range("H136").Value = 1359.2
$ws.Range("J136").Value = 1700
$ws.Range("L136").Value = 5100
$ws.Range("N136").Value = -10200
$ws.Range("H141").Value = 28948.25
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 28948.25
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 28948.25
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -39308.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1100
$ws.Range("J97").Value = 1250
$ws.Range("L97").Value = 3750
$ws.Range("N97").Value = -4742
$ws.Range("H116").Value = 2366.1428
$ws.Range("I116").Value = 1833
$ws.Range("K116").Value = 5499
$ws.Range("M116").Value = -2057
$ws.Range("H120").Value = 10032.444
$ws.Range("I120").Value = 2300
$ws.Range("J120").Value = 10999
$ws.Range("K120").Value = 6900
$ws.Range("L120").Value = 32997
$ws.Range("M120").Value = -2062
$ws.Range("N120").Value = -42673
$ws.Range("H131").Value = 19233692
$ws.Range("J131").Value = 3433.8604
$ws.Range("L131").Value = 10301.5812
$ws.Range("N131").Value = -20381.5812
$ws.Range("H136").Value = 1441
$ws.Range("I136").Value = 965
$ws.Range("J136").Value = 2393
$ws.Range("K136").Value = 2895
$ws.Range("L136").Value = 7179
$ws.Range("M136").Value = 2205
$ws.Range("N136").Value = -17379

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1275.1666
$ws.Range("I22").Value = 750
$ws.Range("J22").Value = 1380.2
$ws.Range("K22").Value = 750
$ws.Range("L22").Value = 1380.2
$ws.Range("M22").Value = -455
$ws.Range("N22").Value = -1970.2
$ws.Range("H27").Value = 1275.1666
$ws.Range("I27").Value = 750
$ws.Range("J27").Value = 1380.2
$ws.Range("K27").Value = 750
$ws.Range("L27").Value = 1380.2
$ws.Range("M27").Value = -643
$ws.Range("N27").Value = -1594.2
$ws.Range("H46").Value = 2683.1667
$ws.Range("I46").Value = 1050
$ws.Range("J46").Value = 3499.75
$ws.Range("K46").Value = 1050
$ws.Range("L46").Value = 3499.75
$ws.Range("M46").Value = -862
$ws.Range("N46").Value = -3875.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H132").Value = 1689.0385
$ws.Range("I132").Value = 1413.7826
$ws.Range("K132").Value = 4241.3478
$ws.Range("M132").Value = -1711.3478
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

